{"js": "// 1. Remove the existing \"_GoBack\" bookmark that wraps \"DONE\" in the\n//    \"zru\u0161i\u0165 administration_view - DONE\" paragraph.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2. Highlight (yellow) the \"Prida\u0165 nov\u00fa rolu ...\" paragraph and append\n//    \" -\" and \" DONE\" (bold) runs, also highlighted.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nlet targetParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].load(\"text\");\n}\nawait context.sync();\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Prida\u0165 nov\u00fa rolu\") !== -1) {\n    targetParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\ntargetParagraph.font.highlightColor = \"Yellow\";\nawait context.sync();\n\nconst dashRun = targetParagraph.insertText(\" -\", \"End\");\ndashRun.font.highlightColor = \"Yellow\";\nawait context.sync();\n\nconst doneRun = targetParagraph.insertText(\" DONE\", \"End\");\ndoneRun.font.highlightColor = \"Yellow\";\ndoneRun.font.bold = true;\nawait context.sync();\n\n// 3. Split the \"dorobi\u0165 funkciu get_user_id...\" run after \"id-\u010d\" and\n//    insert a new \"_GoBack\" bookmark at that point (collapsed, no text\n//    between bookmarkStart/bookmarkEnd).\nconst searchResults = body.search(\"id-\u010d\", { matchCase: true });\nsearchResults.load(\"items\");\nawait context.sync();\n\nconst splitPoint = searchResults.items[0].getRange(\"End\");\nsplitPoint.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Remove the existing \"_GoBack\" bookmark that wraps \"DONE\" in the\n#    \"zru\u0161i\u0165 administration_view - DONE\" paragraph.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# 2. Highlight (yellow) the \"Prida\u0165 nov\u00fa rolu ...\" paragraph and append\n#    \" -\" and \" DONE\" (bold) runs, also highlighted.\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*Prida\u0165 nov\u00fa rolu*\") {\n        $target = $p\n    }\n}\n$target.Range.Font.HighlightColorIndex = 7\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Text = \"pot.\u010dlen\"\n$rng.Find.Execute() | Out-Null\n$rng.Collapse(0)\n\n$rng.InsertAfter(\" -\")\n$rng.Font.HighlightColorIndex = 7\n\n$rng.Collapse(0)\n$rng.InsertAfter(\" DONE\")\n$rng.Font.HighlightColorIndex = 7\n$rng.Font.Bold = 1\n\n# 3. Split the \"dorobi\u0165 funkciu get_user_id...\" run after \"id-\u010d\" and\n#    insert a new \"_GoBack\" bookmark at that point (collapsed, no text\n#    between bookmarkStart/bookmarkEnd).\n$splitRng = $d.Content\n$splitRng.Find.ClearFormatting()\n$splitRng.Find.Text = \"id-\u010d\"\n$splitRng.Find.Execute() | Out-Null\n$splitRng.Collapse(0)\n$d.Bookmarks.Add(\"_GoBack\", $splitRng)\n"}
